$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1863.6666
$ws.Range("I20").Value = 1863.6666
$ws.Range("K20").Value = 1863.6666
$ws.Range("M20").Value = -1633.6666
$ws.Range("H35").Value = 1863.6666
$ws.Range("I35").Value = 1863.6666
$ws.Range("K35").Value = 1863.6666
$ws.Range("M35").Value = -1484.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 651.5
$ws.Range("I12").Value = 1003
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 1003
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = -830
$ws.Range("N12").Value = -646
$ws.Range("H61").Value = 197502.36
$ws.Range("I61").Value = 1273.425
$ws.Range("J61").Value = 911062.0600000001
$ws.Range("K61").Value = 1273.425
$ws.Range("L61").Value = 911062.0600000001
$ws.Range("M61").Value = -1061.425
$ws.Range("N61").Value = -911486.0600000001
$ws.Range("H63").Value = 1795.6
$ws.Range("I63").Value = 1858.25
$ws.Range("K63").Value = 1858.25
$ws.Range("M63").Value = -1172.25
$ws.Range("H66").Value = 1795.6
$ws.Range("I66").Value = 1858.25
$ws.Range("K66").Value = 9291.25
$ws.Range("M66").Value = -5859.25
$ws.Range("H136").Value = 197502.36
$ws.Range("I136").Value = 1273.425
$ws.Range("J136").Value = 911062.0600000001
$ws.Range("K136").Value = 3820.275
$ws.Range("L136").Value = 2733186.18
$ws.Range("M136").Value = -1270.275
$ws.Range("N136").Value = -2738286.18

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3204.4075
$ws.Range("I86").Value = 3010.8572
$ws.Range("J86").Value = 3272.15
$ws.Range("K86").Value = 3010.8572
$ws.Range("L86").Value = 3272.15
$ws.Range("M86").Value = -1887.8572
$ws.Range("N86").Value = -5518.15
$ws.Range("H89").Value = 3204.4075
$ws.Range("I89").Value = 3010.8572
$ws.Range("J89").Value = 3272.15
$ws.Range("K89").Value = 15054.286
$ws.Range("L89").Value = 16360.75
$ws.Range("M89").Value = -9438.286
$ws.Range("N89").Value = -27592.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 33.57895
$ws.Range("I7").Value = 24.785715
$ws.Range("J7").Value = 58.2
$ws.Range("K7").Value = 24.785715
$ws.Range("L7").Value = 58.2
$ws.Range("M7").Value = 88.214285
$ws.Range("N7").Value = -284.2
$ws.Range("H60").Value = 29750
$ws.Range("J60").Value = 29750
$ws.Range("L60").Value = 29750
$ws.Range("N60").Value = -30772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 721.04346
$ws.Range("I34").Value = 365.33334
$ws.Range("J34").Value = 846.58826
$ws.Range("K34").Value = 1096.00002
$ws.Range("L34").Value = 2539.76478
$ws.Range("M34").Value = -1012.00002
$ws.Range("N34").Value = -2707.76478
$ws.Range("H39").Value = 2879.3333
$ws.Range("J39").Value = 3027.8572
$ws.Range("L39").Value = 9083.571599999999
$ws.Range("N39").Value = -9671.571599999999
$ws.Range("H44").Value = 47619132
$ws.Range("I44").Value = 98.666664
$ws.Range("K44").Value = 295.999992
$ws.Range("M44").Value = 102.000008
$ws.Range("H55").Value = 2181.7144
$ws.Range("I55").Value = 516.3333
$ws.Range("J55").Value = 2635.9092
$ws.Range("K55").Value = 1548.9999
$ws.Range("L55").Value = 7907.7276
$ws.Range("M55").Value = -1371.9999
$ws.Range("N55").Value = -8261.7276
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 200
$ws.Range("I60").Value = 200
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 600
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -349
$ws.Range("N60").ClearContents()
$ws.Range("H62").Value = 2765.5
$ws.Range("J62").Value = 2765.5
$ws.Range("L62").Value = 8296.5
$ws.Range("N62").Value = -9668.5
$ws.Range("H63").Value = 1416.8
$ws.Range("I63").Value = 808
$ws.Range("J63").Value = 2330
$ws.Range("K63").Value = 2424
$ws.Range("L63").Value = 6990
$ws.Range("M63").Value = -1675
$ws.Range("N63").Value = -8488
$ws.Range("H64").Value = 33335380
$ws.Range("I64").Value = 1409
$ws.Range("J64").Value = 55558028
$ws.Range("K64").Value = 4227
$ws.Range("L64").Value = 166674084
$ws.Range("M64").Value = -3957
$ws.Range("N64").Value = -166674624
$ws.Range("H65").Value = 2765.5
$ws.Range("J65").Value = 2765.5
$ws.Range("L65").Value = 24889.5
$ws.Range("N65").Value = -31753.5
$ws.Range("H66").Value = 1416.8
$ws.Range("I66").Value = 808
$ws.Range("J66").Value = 2330
$ws.Range("K66").Value = 7272
$ws.Range("L66").Value = 20970
$ws.Range("M66").Value = -3528
$ws.Range("N66").Value = -28458
$ws.Range("H67").Value = 33335380
$ws.Range("I67").Value = 1409
$ws.Range("J67").Value = 55558028
$ws.Range("K67").Value = 4227
$ws.Range("L67").Value = 166674084
$ws.Range("M67").Value = -3291
$ws.Range("N67").Value = -166675956

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13141.75
$ws.Range("I22").Value = 275
$ws.Range("J22").Value = 17430.666
$ws.Range("K22").Value = 275
$ws.Range("L22").Value = 17430.666
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = -18020.666
$ws.Range("H27").Value = 13141.75
$ws.Range("I27").Value = 275
$ws.Range("J27").Value = 17430.666
$ws.Range("K27").Value = 275
$ws.Range("L27").Value = 17430.666
$ws.Range("M27").Value = -168
$ws.Range("N27").Value = -17644.666
$ws.Range("H40").Value = 2570.9285
$ws.Range("I40").Value = 2364.3635
$ws.Range("J40").Value = 3328.3333
$ws.Range("K40").Value = 2364.3635
$ws.Range("L40").Value = 3328.3333
$ws.Range("M40").Value = -2228.3635
$ws.Range("N40").Value = -3600.3333
$ws.Range("H46").Value = 3540.25
$ws.Range("I46").Value = 3720.3333
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 3720.3333
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -3532.3333
$ws.Range("N46").Value = -3376
$ws.Range("H55").Value = 277.26315
$ws.Range("I55").Value = 236.42857
$ws.Range("J55").Value = 301.08334
$ws.Range("K55").Value = 236.42857
$ws.Range("L55").Value = 301.08334
$ws.Range("M55").Value = -63.42857000000001
$ws.Range("N55").Value = -647.08334
$ws.Range("H82").Value = 1325.1052
$ws.Range("I82").Value = 1040.8
$ws.Range("J82").Value = 1641
$ws.Range("K82").Value = 1040.8
$ws.Range("L82").Value = 1641
$ws.Range("M82").Value = -679.8
$ws.Range("N82").Value = -2363
$ws.Range("H85").Value = 1325.1052
$ws.Range("I85").Value = 1040.8
$ws.Range("J85").Value = 1641
$ws.Range("K85").Value = 1040.8
$ws.Range("L85").Value = 1641
$ws.Range("M85").Value = 207.2
$ws.Range("N85").Value = -4137
$ws.Range("H122").Value = 70397.60000000001
$ws.Range("I122").Value = 145772
$ws.Range("J122").Value = 4445
$ws.Range("K122").Value = 13335
$ws.Range("L122").Value = 13335
$ws.Range("M122").Value = -434866
$ws.Range("N122").Value = -18235

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 675
$ws.Range("I13").Value = 450
$ws.Range("J13").Value = 900
$ws.Range("K13").Value = 450
$ws.Range("L13").Value = 900
$ws.Range("M13").Value = -310
$ws.Range("N13").Value = -1180
$ws.Range("H82").Value = 46500
$ws.Range("J82").Value = 46500
$ws.Range("L82").Value = 46500
$ws.Range("N82").Value = -47266
$ws.Range("H85").Value = 46500
$ws.Range("J85").Value = 46500
$ws.Range("L85").Value = 46500
$ws.Range("N85").Value = -49152
